$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily log entry (2026/01/05, Mon, 13:00, 15) was inserted as a new
# row 563, pushing the existing rows 563:604 down to 564:605 and extending
# the used range from A1:D604 to A1:D605.
$ws.Rows.Item(563).Insert()

# Write the date as text (matching the existing "YYYY/MM/DD" text entries
# in column A) rather than letting Excel auto-convert it to a date serial.
$ws.Cells.Item(563, 1).NumberFormat = "@"
$ws.Cells.Item(563, 1).Value = "2026/01/05"
$ws.Cells.Item(563, 1).Style = "Normal"

$ws.Cells.Item(563, 2).Value = "月"
$ws.Cells.Item(563, 3).Value = 13
$ws.Cells.Item(563, 4).Value = 15
